$d = $word.ActiveDocument

# --- Edit 1: paragraph 16 ("Nu-uh, dead. Dead. Dead.") becomes a split reply ---
$p16 = $d.Paragraphs.Item(16)
$r16 = $p16.Range
$xmlFrag1 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t>“Nu-uh,</w:t></w:r><w:r><w:t xml:space="preserve"> you’re dead,” he replied, giving her a few more stabs for emphasis.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r16.InsertXML($xmlFrag1)

# --- Edit 2: paragraphs 17-20 (Moonsong's soliloquy .. trailing blank paragraphs)
# get replaced with the new housebreak scene + relocated bookmark ---
$p17 = $d.Paragraphs.Item(17)
$p20 = $d.Paragraphs.Item(20)
$r2 = $d.Range($p17.Range.Start, $p20.Range.End)
$xmlFrag2 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:ind w:firstLine="720"/></w:pPr><w:r><w:t xml:space="preserve">Moonsong’s soliloquy was never meant to be. She took the last spear with a flourish, rolled over, and died with Bloodfang’s name on her lips. Squinting, she saw Percy open </w:t></w:r><w:r><w:t>an eye</w:t></w:r><w:r><w:t xml:space="preserve"> and smile at her blood offering to the performing arts.</w:t></w:r><w:r><w:t xml:space="preserve"> By way of response, Vera let out another death choke, convulsed, and died again.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:i/></w:rPr></w:pPr></w:p><w:p><w:r><w:rPr><w:i/></w:rPr><w:tab/></w:r><w:r><w:t>The hacksaw was perhaps excessive, but after her parents had padlocked the ground-level windows, she’d not been left with many options. They would be cross with her, but the knot in her stomach didn’</w:t></w:r><w:r><w:t>t care anymore. The rabbit tried to ignore the nervous prickling under her fur as she pulled it back and forth</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> much too loud for her liking. Perhaps the rattling bone jewelry had been poorly advised too, but she felt better with it on. With any luck, the next iteration would be with bones from her own kills. </w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>The padlock clattered t</w:t></w:r><w:r><w:t xml:space="preserve">o the ground. Vera </w:t></w:r><w:r><w:t>stacked the evidence</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t xml:space="preserve"> threw a look over her shoulder,</w:t></w:r><w:r><w:t xml:space="preserve"> and threw open the window. </w:t></w:r><w:r><w:t xml:space="preserve">She drank in the night air, and already her racing heart began to slow. </w:t></w:r><w:r><w:t xml:space="preserve">Vera’s feet were already half out the window before she remembered to unsling her pack, which was so full to bursting that it would be </w:t></w:r><w:r><w:t>im</w:t></w:r><w:r><w:t>possible to fit through still attached to he</w:t></w:r><w:r><w:t>r</w:t></w:r><w:r><w:t>. Her hands fumbled with the straps</w:t></w:r><w:r><w:t xml:space="preserve"> as soon as it was out</w:t></w:r><w:r><w:t xml:space="preserve">, throwing the pack, bow, </w:t></w:r><w:r><w:t>and makeshift quiver back into place.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>She raced for the tree line, thankful that her warren was uncharacteristically close to the woods, uncharacteristically close to the predators and dangerous creatures of the night that Vera, already a night owl, had been told to fear. She pushed forward, putting as much distance between herself and the elaborate earthen tomb as she could. Much more so than any of the scents of the warren, it was the crisp scent of leaves and the touch of cool air that made Vera feel at home.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Once she was certain she was out of the range of prying eyes, the rabbit unslung her backpack and made to</w:t></w:r><w:r><w:t xml:space="preserve"> set up camp. There was a certain excitement to it, knowing that o</w:t></w:r><w:r><w:t xml:space="preserve">ut here, she really was exposed, that life was not docile, </w:t></w:r><w:r><w:t>predictable, and frustratingly safe inside four packed earth walls</w:t></w:r><w:r><w:t>.</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t xml:space="preserve">She </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t>conjured a fire through the witchcraft of flint and steel, then laid down, pulling out a worn leather notebook from her bag. Her bow came off, never more than an arm’s reach away.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:tab/></w:r><w:r><w:t>The thought that whatever encountered her probably had more to fear from her than her from them brought a grin to her face, and for a moment, she was living her fantasy, she was Moonsong.</w:t></w:r><w:r><w:t xml:space="preserve"> She flipped through the book – illustrations and pages of notes and unfinished stories – until she got to a half-finished page near the end. She sighed wistfully, uncorking an ink flask and preparing her pen, when something caught her eye.</w:t></w:r></w:p><w:p><w:r><w:tab/><w:t>Golden eyes shone back at her.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$r2.InsertXML($xmlFrag2)
